# Applies the "Actualizacion automatica" data refresh described by the diff.
# The workbook has three sheets:
#   1. VENTAS POR GRUPO     (raw per-client-per-product-group sales)
#   2. VENTA MENSUAL        (per-client monthly sales totals + budget)
#   3. CUMPLIMIENTO MENSUAL (per-product-group budget vs. sales summary)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("M4").Value = 1924.35

$ws1.Range("C12").Value = 518.4
$ws1.Range("L12").Value = 1089.41

$ws1.Range("M13").Value = 2683.55

$ws1.Range("C14").Value = 414.72
$ws1.Range("M14").Value = -636.72

$ws1.Range("M15").Value = 2762.6

$ws1.Range("D26").Value = 933.12
$ws1.Range("M26").Value = 1280.2

$ws1.Range("L29").Value = 855.36
$ws1.Range("M29").Value = 1382.34

$ws1.Range("C54").Value = "2 de 52"
$ws1.Range("L54").Value = "8 de 52"
$ws1.Range("M54").Value = "13 de 52"

# ---------------------------------------------------------------------
# Sheet 2: VENTA MENSUAL
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("F4").Value = 4227.03
$ws2.Range("F12").Value = 1607.81
$ws2.Range("F13").Value = 2683.55
$ws2.Range("F14").Value = 1305.91
$ws2.Range("F15").Value = 2762.6
$ws2.Range("F26").Value = 2213.32
$ws2.Range("F29").Value = 2237.7
$ws2.Range("F58").Value = 44033.81

# ---------------------------------------------------------------------
# Sheet 3: CUMPLIMIENTO MENSUAL
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# row 2: 240X120 PORCELANATO
$ws3.Range("D2").Value = 933.12
$ws3.Range("E2").Value = 5264.46402943659
$ws3.Range("F2").Value = 0.1505618956625632

# row 3: 240X80 PORCELANATO
$ws3.Range("D3").Value = 2697.4
$ws3.Range("E3").Value = 14971.7470988183
$ws3.Range("F3").Value = 0.1526615849035747

# row 11: PIEDRA SINTERIZADA
$ws3.Range("D11").Value = 9691.5
$ws3.Range("E11").Value = 8139.9143984654
$ws3.Range("F11").Value = 0.5435070815713904

# row 12: PORCELANATO
$ws3.Range("D12").Value = 23143.39
$ws3.Range("E12").Value = 38720.3303947566
$ws3.Range("F12").Value = 0.3741027835429304

# row 15: TOTAL
$ws3.Range("D15").Value = 42883.05
$ws3.Range("E15").Value = 79171.78551083435
$ws3.Range("F15").Value = 0.3513424914344621
